# Commit: "Add unique VINs to Each CA Select Test PT1"
#
# Rows 2-5 on Sheet1 all shared the same placeholder VIN value
# "1FDEU15H&K" in column A. Give each of those rows a unique VIN so the
# test data no longer collides - replace it with a brand-new VIN value
# "AAANK3CC&F".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AAANK3CC&F"
$ws.Range("A3").Value = "AAANK3CC&F"
$ws.Range("A4").Value = "AAANK3CC&F"
$ws.Range("A5").Value = "AAANK3CC&F"

# Leave the cursor where the author left it when saving the file.
[void]$ws.Range("B12").Select()
